$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing product rows ---------------------------------------
# Row 2: Chicken -> Pasta sauce, quantity bumped to 15
$ws.Range("A2").Value = "PASTA SAUCE-350 GRM"
$ws.Range("B2").Value = 15.0

# Row 3: Beef -> Fish sauce, unit price corrected
$ws.Range("A3").Value = "FISH SAUCE-600 ML"
$ws.Range("C3").Value = 600.0

# The Remark cells were mistakenly placed out in column G; bring them back
# under the Remark header (column E) next to the rest of the row data.
$ws.Range("G2").Cut($ws.Range("E2"))
$ws.Range("G2").Clear()
$ws.Range("G3").Cut($ws.Range("E3"))
$ws.Range("G3").Clear()

# --- Add newly confirmed inquiry rows ------------------------------------
# Use the existing row formatting as the template for the two new rows.
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

$ws.Range("A4").Value = "PRIMA KOTTU MEE - 80GRM"
$ws.Range("B4").Value = 30.0
$ws.Range("C4").Value = 140.0
$ws.Range("D4").Value = "Available"

$ws.Range("A5").Value = "PIZZA BASE-5 PCS"
$ws.Range("B5").Value = 30.0
$ws.Range("C5").Value = 850.0
$ws.Range("D5").Value = "Available"

# Row 5's remark holds real confirmation text, so give it the same style
# used by the other descriptive (Status/Remark) cells before filling it in.
$ws.Range("D3").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = "There is a 10% discount for 20+"

# --- Resize columns to fit the new, longer product/remark text ----------
# (ColumnWidth is snapped to the host's pixel grid on write, so the inputs
# below are chosen to land on the grid point nearest the intended widths:
# ~26.57, ~18.86 and ~42.71 "characters" respectively.)
$ws.Columns.Item(1).ColumnWidth = 25.665
$ws.Columns.Item(4).ColumnWidth = 18.0
$ws.Columns.Item(5).ColumnWidth = 41.83
